# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates column G ("K") on the active worksheet for rows 2-74 with the
# freshly-calculated strike-count values (s_vals), replacing the stale
# "Strike#" figures that were previously stored there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0;  3  = 0;  4  = 2;  5  = 1;  6  = 1;  7  = 1;  8  = 1;  9  = 2;  10 = 3;
    11 = 0;  12 = 1;  13 = 2;  14 = 0;  15 = 2;  16 = 1;  17 = 1;  18 = 1;  19 = 0;
    20 = 1;  21 = 1;  22 = 1;  23 = 1;  24 = 1;  25 = 2;  26 = 2;  27 = 0;  28 = 3;
    29 = 0;  30 = 2;  31 = 2;  32 = 2;  33 = 1;  34 = 0;  35 = 2;  36 = 2;  37 = 1;
    38 = 1;  39 = 1;  40 = 2;  41 = 0;  42 = 2;  43 = 2;  44 = 2;  45 = 1;  46 = 1;
    47 = 0;  48 = 0;  49 = 0;  50 = 1;  51 = 1;  52 = 1;  53 = 0;  54 = 1;  55 = 1;
    56 = 1;  57 = 1;  58 = 2;  59 = 2;  60 = 3;  61 = 2;  62 = 1;  63 = 1;  64 = 2;
    65 = 1;  66 = 2;  67 = 1;  68 = 1;  69 = 0;  70 = 2;  71 = 1;  72 = 1;  73 = 3;
    74 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
